# Refresh the "cryptos" price/volume snapshot (GitHub Actions style data refresh).
# Column D (Price) and E (Volume(1h)) are stored as plain text in this sheet, so for
# every Price cell we force NumberFormat "@" (Text) before assigning the new value -
# otherwise Excel would auto-parse numeric-looking strings like "1.00" or "64.237.70"
# into real numbers and lose the original text formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.237.70'
$ws.Range("E2").Value = '  -2.66%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.176.78'
$ws.Range("E3").Value = '  -7.85%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.16'
$ws.Range("E5").Value = '  -4.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.73'
$ws.Range("E6").Value = '  -1.85%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.600'
$ws.Range("E8").Value = '  -0.23%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.175.64'
$ws.Range("E9").Value = '  -7.82%  '

$ws.Range("E10").Value = '  -6.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.60'
$ws.Range("E11").Value = '  -4.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.395'
$ws.Range("E12").Value = '  -3.75%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.727.39'
$ws.Range("E13").Value = '  -7.86%  '

$ws.Range("E14").Value = '  -0.03%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.47'
$ws.Range("E15").Value = '  -5.52%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.320.46'
$ws.Range("E16").Value = '  -2.60%  '

$ws.Range("E17").Value = '  -5.57%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.160.32'
$ws.Range("E18").Value = '  -8.31%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.64'
$ws.Range("E19").Value = '  -4.99%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.97'
$ws.Range("E20").Value = '  -6.28%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '351.83'
$ws.Range("E21").Value = '  -4.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.18'
$ws.Range("E22").Value = '  -6.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.00'
$ws.Range("E24").Value = '  -5.16%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.501'
$ws.Range("E25").Value = '  -6.63%  '

$ws.Range("E26").Value = '  -3.81%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.38'
$ws.Range("E27").Value = '  -3.74%  '

$ws.Range("E28").Value = '  -0.79%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.08%  '

$ws.Range("E31").Value = '  -2.74%  '

$ws.Range("E32").Value = '  -4.75%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '22.08'
$ws.Range("E33").Value = '  -7.25%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.58'
$ws.Range("E34").Value = '  -6.32%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.18'
$ws.Range("E35").Value = '  -8.71%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '158.14'
$ws.Range("E36").Value = '  -1.59%  '

$ws.Range("E37").Value = '  -6.89%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.805'
$ws.Range("E38").Value = '  -8.54%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.06'
$ws.Range("E39").Value = '  -9.72%  '

$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.48'
$ws.Range("E40").Value = '  -6.28%  '

$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.67'
$ws.Range("E41").Value = '  -5.09%  '

$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.644.25'
$ws.Range("E42").Value = '  -4.31%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.99'
$ws.Range("E43").Value = '  -6.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.13'
$ws.Range("E44").Value = '  -7.38%  '

$ws.Range("E45").Value = '  -4.23%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.60'
$ws.Range("E46").Value = '  -3.86%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '321.74'
$ws.Range("E47").Value = '  -0.79%  '

$ws.Range("E48").Value = '  -4.35%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0269'
$ws.Range("E49").Value = '  -7.09%  '

$ws.Range("E50").Value = '  -0.32%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  +0.03%  '
